# Insert a new weekly record at row 427 for "Hortaliza, Vega Monumental
# Concepción - Coliflor". This pushes the existing rows 427:505 down to
# 428:506 (dimension grows from A1:R505 to A1:R506) and fills the newly
# inserted row 427 with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 427 downward (this also grows the sheet dimension automatically).
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new data point.
$ws.Range("A427").Value = 11
$ws.Range("B427").Value = "Vega Monumental Concepción"
$ws.Range("C427").Value = "Bíobío"
$ws.Range("D427").Value = 45244
$ws.Range("E427").Value = 8
$ws.Range("F427").Value = 100112008
$ws.Range("G427").Value = "Coliflor"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 1500
$ws.Range("K427").Value = 1000
$ws.Range("L427").Value = 1000
$ws.Range("M427").Value = 1000
$ws.Range("N427").Value = "$/unidad"
$ws.Range("O427").Value = "Región Metropolitana"
$ws.Range("P427").Value = 1000
$ws.Range("Q427").Value = 1
$ws.Range("R427").Value = "Hortaliza"

Write-Host "Inserted new row 427 and shifted subsequent rows down; dimension now A1:R506"
